$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 (shifts existing rows 14-18 down to 15-19)
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the Chloride [as Cl] / 00940 parameter code entry
$ws.Range("A14").Value = "00940"
$ws.Range("B14").Value = "Chloride [as Cl]"
$ws.Range("C14").Value = "Chlorides"
$ws.Range("D14").Value = "mg/L"

# Update the selection to match the post-edit state (A20 was selected after the new last row 19)
$ws.Range("A20").Select()
